$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.925.36'
$ws.Range("E2").Value = '  +2.12%  '
$ws.Range("D3").Value = '3.474.05'
$ws.Range("E3").Value = '  +2.33%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.93'
$ws.Range("E5").Value = '  +0.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.18'
$ws.Range("E6").Value = '  +3.98%  '
$ws.Range("D7").Value = '3.474.30'
$ws.Range("E7").Value = '  +2.33%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.481'
$ws.Range("E9").Value = '  +1.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.69'
$ws.Range("E10").Value = '  +0.83%  '
$ws.Range("E11").Value = '  +1.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.402'
$ws.Range("E12").Value = '  +4.10%  '
$ws.Range("D13").Value = '4.065.86'
$ws.Range("E13").Value = '  +2.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.81'
$ws.Range("E14").Value = '  +6.63%  '
$ws.Range("E15").Value = '  +2.82%  '
$ws.Range("D16").Value = '3.473.61'
$ws.Range("E16").Value = '  +2.21%  '
$ws.Range("E17").Value = '  +0.56%  '
$ws.Range("D18").Value = '62.900.38'
$ws.Range("E18").Value = '  +2.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.32'
$ws.Range("E19").Value = '  +3.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.38'
$ws.Range("E20").Value = '  +5.55%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.22'
$ws.Range("E21").Value = '  +1.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '389.41'
$ws.Range("E22").Value = '  +0.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.558'
$ws.Range("E23").Value = '  +1.90%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.68'
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("D26").Value = '3.604.11'
$ws.Range("E26").Value = '  +2.07%  '
$ws.Range("E27").Value = '  +1.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.178'
$ws.Range("E28").Value = '  +0.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.58'
$ws.Range("E29").Value = '  +2.60%  '
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.15'
$ws.Range("E31").Value = '  +2.32%  '
$ws.Range("E32").Value = '  -0.58%  '
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("E34").Value = '  -3.24%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.63'
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.26'
$ws.Range("E36").Value = '  +3.68%  '
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.06'
$ws.Range("E37").Value = '  +2.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '31.95'
$ws.Range("E38").Value = '  +18.56%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '170.72'
$ws.Range("E39").Value = '  +1.41%  '
$ws.Range("E40").Value = '  +6.41%  '
$ws.Range("D41").Value = '3.509.35'
$ws.Range("E41").Value = '  +2.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0755'
$ws.Range("E42").Value = '  -1.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.800'
$ws.Range("E43").Value = '  +2.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.41'
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.46'
$ws.Range("E45").Value = '  +0.39%  '
$ws.Range("E46").Value = '  +2.81%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.21'
$ws.Range("E47").Value = '  +3.91%  '
$ws.Range("D48").Value = '2.611.98'
$ws.Range("E48").Value = '  +5.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.28'
$ws.Range("E49").Value = '  +12.66%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.94'
$ws.Range("E50").Value = '  +1.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.72'
$ws.Range("E51").Value = '  +1.01%  '
